$wb = $excel.ActiveWorkbook

# --- Sheet: accuracy ---
$ws = $wb.Worksheets.Item("accuracy")
$ws.Range("B2").Value = [double]"1.618173175280639e-08"
$ws.Range("B3").Value = [double]"0.0001511117315160743"
$ws.Range("B4").Value = [double]"3.246171205575621e-05"
$ws.Range("B5").Value = [double]"0.09124652022007672"
$ws.Range("B6").Value = [double]"0.1001309985781935"
$ws.Range("B7").Value = [double]"0.01176300432789639"
$ws.Range("B8").Value = [double]"0.1183057268373853"
$ws.Range("B9").Value = [double]"0.5770324017264793"
$ws.Range("B10").Value = [double]"9.993281145492068e-05"
$ws.Range("B11").Value = [double]"3.863341453943772e-05"
$ws.Range("B12").Value = [double]"9.620163534427312e-05"
$ws.Range("B13").Value = [double]"0.02182190128758563"
$ws.Range("B14").Value = [double]"0.001254375329880237"
$ws.Range("B15").Value = [double]"1.005025225339669e-08"

# --- Sheet: sensitivity ---
$ws = $wb.Worksheets.Item("sensitivity")
$ws.Range("B2").Value = [double]"0.0001602085555052488"
$ws.Range("B3").Value = [double]"3.471063405753461e-13"
$ws.Range("B4").Value = [double]"6.203442143418255e-11"
$ws.Range("B5").Value = [double]"5.691335972560574e-08"
$ws.Range("B6").Value = [double]"9.260710188887272e-05"
$ws.Range("B7").Value = [double]"7.061502350243177e-12"
$ws.Range("B8").Value = [double]"2.358169682807503e-07"
$ws.Range("B9").Value = [double]"1.002723295488438e-08"
$ws.Range("B10").Value = [double]"0.4894504311402167"
$ws.Range("B11").Value = [double]"3.555597563440914e-08"
$ws.Range("B12").Value = [double]"2.684135272093207e-13"
$ws.Range("B13").Value = [double]"6.049066718272365e-10"
$ws.Range("B14").Value = [double]"4.955163747957975e-13"
$ws.Range("B15").Value = [double]"2.09523121127742e-16"

# --- Sheet: specificity ---
$ws = $wb.Worksheets.Item("specificity")
$ws.Range("B2").Value = [double]"6.141723331800252e-08"
$ws.Range("B3").Value = [double]"1.333170065933796e-10"
$ws.Range("B4").Value = [double]"2.710274073500794e-10"
$ws.Range("B5").Value = [double]"9.993419997549538e-06"
$ws.Range("B6").Value = [double]"4.084616518510084e-07"
$ws.Range("B7").Value = [double]"7.966900069562025e-08"
$ws.Range("B8").Value = [double]"0.0002574540638237679"
$ws.Range("B9").Value = [double]"2.10491016137352e-06"
$ws.Range("B10").Value = [double]"1.379820435783134e-07"
$ws.Range("B11").Value = [double]"0.003634820685442738"
$ws.Range("B12").Value = [double]"2.010515489085699e-12"
$ws.Range("B13").Value = [double]"2.303678704473348e-08"
$ws.Range("B14").Value = [double]"1.328322046511736e-09"
$ws.Range("B15").Value = [double]"1.733353365758601e-14"

# --- Sheet: time ---
$ws = $wb.Worksheets.Item("time")
$ws.Range("B2:B15").Value = [double]"6.192194413718758e-34"
